# KK-SEM Ultrahang HW BOM — reshuffle capacitor designators, add C12 and C26
# (new MIC 100nF cap + new 10u cap), renumber the 100n/10u/1u/10n rows.
# Matches commit "Added MIC 100nF, U2D ref".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capacitor table (rows 13-18) ---------------------------------------

# Row 13: 1n-50V, C0G 1% group gains C12 (loses C8, which moves to 100n group)
$ws.Range("D13").Value = "C1, C2, C3, C4, C5, C6, C9, C10, C11, C12"

# Row 14: 100n-50V group now has 8 parts (was 7) -> gains C8, C14, C17, C21, C23, C25
$ws.Range("A14").Formula = '=$M$1*8'
$ws.Range("D14").Value = "C7, C8, C13, C14, C17, C21, C23, C25"

# Row 15: was 10n-50V/C17, now 10u-16V group (C15, C22, C24, and new C26), X5R
$ws.Range("B15").Value = "10u-16V"
$ws.Range("D15").Value = "C15, C22, C24, C26"
$ws.Range("E15").Value = "X5R"

# Row 16: was 10u-16V group, now 1u-16V/C16, X7R
$ws.Range("B16").Value = "1u-16V"
$ws.Range("D16").Value = "C16"
$ws.Range("E16").Value = "X7R"

# Row 17: was 1u-16V/C15, now 10n-50V/C18
$ws.Range("B17").Value = "10n-50V"
$ws.Range("D17").Value = "C18"

# Row 18: 22p-50V group now C19, C20 (was C18, C19)
$ws.Range("D18").Value = "C19, C20"

# --- Column D best-fit width grew by ~1 char due to the longer D13 text --
$ws.Columns.Item(4).ColumnWidth = 31.6

# --- Selection / viewport, matches the author's final cursor position ---
$ws.Range("D18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
